$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 42-45: four "Tretåig hackspett" observation records. The species /
# location-description fields are identical across the four rows and stay
# put; only the per-record Id (A), coordinates (Q, R) and public comment
# (AC) are cyclically rotated - row 45's old values wrap around to row 42.
# ---------------------------------------------------------------------------

function Get-Rec($row) {
    return [PSCustomObject]@{
        A  = $ws.Range("A$row").Value2
        Q  = $ws.Range("Q$row").Value2
        R  = $ws.Range("R$row").Value2
        AC = $ws.Range("AC$row").Value2
    }
}

function Set-Rec($row, $vals) {
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("AC$row").Value = $vals.AC
}

$rec42 = Get-Rec 42
$rec43 = Get-Rec 43
$rec44 = Get-Rec 44
$rec45 = Get-Rec 45

Set-Rec 42 $rec43
Set-Rec 43 $rec44
Set-Rec 44 $rec45
Set-Rec 45 $rec42

# ---------------------------------------------------------------------------
# Rows 59-61: entire observation records shift down by one row - row 59's
# content moves to row 60, row 60's moves to row 61, and row 61's wraps
# around to row 59. The fields that actually differ between the three rows
# are the Id (A), taxon sort (B), TaxonId (E), species name (F), scientific
# name (G), author (H), coordinates (Q, R) and public comment (AC); every
# other column already holds the same value in all three rows.
# ---------------------------------------------------------------------------

function Get-FullRec($row) {
    return [PSCustomObject]@{
        A  = $ws.Range("A$row").Value2
        B  = $ws.Range("B$row").Value2
        E  = $ws.Range("E$row").Value2
        F  = $ws.Range("F$row").Value2
        G  = $ws.Range("G$row").Value2
        H  = $ws.Range("H$row").Value2
        Q  = $ws.Range("Q$row").Value2
        R  = $ws.Range("R$row").Value2
        AC = $ws.Range("AC$row").Value2
    }
}

function Set-FullRec($row, $vals) {
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    if ($null -eq $vals.AC -or $vals.AC -eq "") {
        $ws.Range("AC$row").ClearContents()
    } else {
        $ws.Range("AC$row").Value = $vals.AC
    }
}

$rec59 = Get-FullRec 59
$rec60 = Get-FullRec 60
$rec61 = Get-FullRec 61

Set-FullRec 59 $rec61
Set-FullRec 60 $rec59
Set-FullRec 61 $rec60
